$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.62
$ws.Range("G2").Value = 3.35
$ws.Range("H2").Value = 2.38
$ws.Range("I2").Value = 3.3
$ws.Range("J2").Value = 2.96
$ws.Range("K2").Value = 3.85
$ws.Range("L2").Value = 1.32
$ws.Range("N2").Value = 1.9
$ws.Range("P2").Value = 1.76
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 1.08
$ws.Range("S2").Value = 1.75
$ws.Range("V2").Value = 1.43
$ws.Range("W2").Value = 1.42
$ws.Range("F3").Value = 1.75
$ws.Range("G3").Value = 2.1
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 980
$ws.Range("K3").Value = 5.6
$ws.Range("N3").Value = 1.66
$ws.Range("P3").Value = 1.66
$ws.Range("R3").Value = 1.1
$ws.Range("T3").Value = 1.03
$ws.Range("U3").Value = 1.03
$ws.Range("V3").Value = 1.12
$ws.Range("W3").Value = 1.91
$ws.Range("L4").Value = 1.37
$ws.Range("N4").Value = 3.85
$ws.Range("R4").Value = 1.36
$ws.Range("X4").Value = 980
$ws.Range("Y4").Value = 980
$ws.Range("AD4").Value = 980
$ws.Range("AF4").Value = 9.800000000000001
$ws.Range("AH4").Value = 980
$ws.Range("AL4").Value = 980
$ws.Range("G5").Value = 2.82
$ws.Range("Q5").Value = 1.77
$ws.Range("T5").Value = 1.03
$ws.Range("U5").Value = 1.03
$ws.Range("W5").Value = 1.55
$ws.Range("F6").Value = 4
$ws.Range("H6").Value = 1.69
$ws.Range("K6").Value = 4.1
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 1.11
$ws.Range("P6").Value = 1.78
$ws.Range("Q6").Value = 1.92
$ws.Range("R6").Value = 1.25
$ws.Range("T6").Value = 1.03
$ws.Range("U6").Value = 1.03
$ws.Range("X7").Value = 980
$ws.Range("Y7").Value = 980
$ws.Range("Z7").Value = 980
$ws.Range("AB7").Value = 980
$ws.Range("AC7").Value = 980
$ws.Range("AD7").Value = 980
$ws.Range("AF7").Value = 980
$ws.Range("AG7").Value = 980
$ws.Range("AH7").Value = 980
$ws.Range("AJ7").Value = 980
$ws.Range("AK7").Value = 980
$ws.Range("AL7").Value = 980
$ws.Range("AN7").Value = 980
$ws.Range("K8").Value = 5.6
$ws.Range("L8").Value = 1.01
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 1.46
$ws.Range("O8").Value = 1.01
$ws.Range("P8").Value = 1.46
$ws.Range("Q8").Value = 2.32
$ws.Range("R8").Value = 1.09
$ws.Range("S8").Value = 2.32
$ws.Range("T8").Value = 1.03
$ws.Range("U8").Value = 1.03
$ws.Range("V8").Value = 1.1
$ws.Range("W8").Value = 1.98
$ws.Range("X8").Value = 1000
$ws.Range("Y8").Value = 1000
$ws.Range("Z8").Value = 1000
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 1000
$ws.Range("AC8").Value = 1000
$ws.Range("AD8").Value = 1000
$ws.Range("AE8").Value = 1000
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 1000
$ws.Range("AH8").Value = 1000
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 1000
$ws.Range("AK8").Value = 1000
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 1000
$ws.Range("AO8").Value = 1000
$ws.Range("G9").Value = 3.35
$ws.Range("H9").Value = 2.14
$ws.Range("L9").Value = 1.01
$ws.Range("M9").Value = 1.01
$ws.Range("N9").Value = 1.79
$ws.Range("O9").Value = 1.01
$ws.Range("R9").Value = 1.23
$ws.Range("S9").Value = 2.96
$ws.Range("T9").Value = 1.56
$ws.Range("U9").Value = 1.76
$ws.Range("V9").Value = 1.42
$ws.Range("W9").Value = 1.42
$ws.Range("X9").Value = 980
$ws.Range("Y9").Value = 1000
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 1000
$ws.Range("AC9").Value = 1000
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 1000
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 1000
$ws.Range("I10").Value = 4.6
$ws.Range("K10").Value = 3.05
$ws.Range("L10").Value = 1.01
$ws.Range("M10").Value = 1.14
$ws.Range("N10").Value = 2.24
$ws.Range("O10").Value = 1.7
$ws.Range("R10").Value = 1.13
$ws.Range("S10").Value = 7.4
$ws.Range("T10").Value = 2.48
$ws.Range("U10").Value = 1.58
$ws.Range("V10").Value = 1.27
$ws.Range("W10").Value = 1.79
$ws.Range("X10").Value = 7
$ws.Range("Y10").Value = 10
$ws.Range("Z10").Value = 30
$ws.Range("AA10").Value = 130
$ws.Range("AB10").Value = 6.2
$ws.Range("AC10").Value = 7.6
$ws.Range("AD10").Value = 23
$ws.Range("AE10").Value = 100
$ws.Range("AF10").Value = 11.5
$ws.Range("AG10").Value = 13
$ws.Range("AH10").Value = 34
$ws.Range("AI10").Value = 180
$ws.Range("AJ10").Value = 32
$ws.Range("AK10").Value = 38
$ws.Range("AL10").Value = 90
$ws.Range("AM10").Value = 380
$ws.Range("AN10").Value = 42
$ws.Range("AO10").Value = 190
$ws.Range("G11").Value = 2.96
$ws.Range("L11").Value = 1.57
$ws.Range("S11").Value = 5.2
$ws.Range("W11").Value = 1.51
$ws.Range("I12").Value = 3.35
$ws.Range("J12").Value = 3.05
$ws.Range("L12").Value = 1.01
$ws.Range("M12").Value = 1.01
$ws.Range("N12").Value = 1.58
$ws.Range("O12").Value = 1.01
$ws.Range("Q12").Value = 2.16
$ws.Range("R12").Value = 1.19
$ws.Range("S12").Value = 3.6
$ws.Range("T12").Value = 1.03
$ws.Range("U12").Value = 1.03
$ws.Range("V12").Value = 1.42
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 1000
$ws.Range("Y12").Value = 1000
$ws.Range("Z12").Value = 1000
$ws.Range("AA12").Value = 1000
$ws.Range("AB12").Value = 1000
$ws.Range("AC12").Value = 1000
$ws.Range("AD12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 1000
$ws.Range("AG12").Value = 1000
$ws.Range("AH12").Value = 1000
$ws.Range("AI12").Value = 1000
$ws.Range("AJ12").Value = 1000
$ws.Range("AK12").Value = 1000
$ws.Range("AL12").Value = 1000
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 1000
$ws.Range("AO12").Value = 1000
$ws.Range("G13").Value = 3.15
$ws.Range("K13").Value = 980
$ws.Range("L13").Value = 1.01
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 1.4
$ws.Range("O13").Value = 1.01
$ws.Range("P13").Value = 1.4
$ws.Range("R13").Value = 1.09
$ws.Range("S13").Value = 2.38
$ws.Range("T13").Value = 1.03
$ws.Range("U13").Value = 1.03
$ws.Range("V13").Value = 1.33
$ws.Range("W13").Value = 1.46
$ws.Range("X13").Value = 1000
$ws.Range("Y13").Value = 1000
$ws.Range("Z13").Value = 1000
$ws.Range("AA13").Value = 1000
$ws.Range("AB13").Value = 1000
$ws.Range("AC13").Value = 1000
$ws.Range("AD13").Value = 1000
$ws.Range("AE13").Value = 1000
$ws.Range("AF13").Value = 1000
$ws.Range("AG13").Value = 1000
$ws.Range("AH13").Value = 1000
$ws.Range("AI13").Value = 1000
$ws.Range("AJ13").Value = 1000
$ws.Range("AK13").Value = 1000
$ws.Range("AL13").Value = 1000
$ws.Range("AM13").Value = 1000
$ws.Range("AN13").Value = 1000
$ws.Range("AO13").Value = 1000
$ws.Range("F14").Value = 1.92
$ws.Range("G14").Value = 2.62
$ws.Range("H14").Value = 3.1
$ws.Range("K14").Value = 6.4
$ws.Range("L14").Value = 1.01
$ws.Range("M14").Value = 1.01
$ws.Range("N14").Value = 1.63
$ws.Range("O14").Value = 1.01
$ws.Range("Q14").Value = 1.93
$ws.Range("R14").Value = 1.08
$ws.Range("S14").Value = 1.93
$ws.Range("T14").Value = 1.01
$ws.Range("U14").Value = 1.01
$ws.Range("V14").Value = 1.28
$ws.Range("W14").Value = 1.61
$ws.Range("X14").Value = 1000
$ws.Range("Y14").Value = 1000
$ws.Range("Z14").Value = 1000
$ws.Range("AA14").Value = 1000
$ws.Range("AB14").Value = 1000
$ws.Range("AC14").Value = 1000
$ws.Range("AD14").Value = 1000
$ws.Range("AE14").Value = 1000
$ws.Range("AF14").Value = 1000
$ws.Range("AG14").Value = 1000
$ws.Range("AH14").Value = 1000
$ws.Range("AI14").Value = 1000
$ws.Range("AJ14").Value = 1000
$ws.Range("AK14").Value = 1000
$ws.Range("AL14").Value = 1000
$ws.Range("AM14").Value = 1000
$ws.Range("AN14").Value = 1000
$ws.Range("AO14").Value = 1000
$ws.Range("H15").Value = 2.56
$ws.Range("K15").Value = 5.1
$ws.Range("L15").Value = 1.01
$ws.Range("M15").Value = 1.01
$ws.Range("N15").Value = 1.53
$ws.Range("O15").Value = 1.01
$ws.Range("R15").Value = 1.08
$ws.Range("S15").Value = 2.06
$ws.Range("T15").Value = 1.01
$ws.Range("U15").Value = 1.01
$ws.Range("V15").Value = 1.28
$ws.Range("W15").Value = 1.42
$ws.Range("X15").Value = 1000
$ws.Range("Y15").Value = 1000
$ws.Range("Z15").Value = 1000
$ws.Range("AA15").Value = 1000
$ws.Range("AB15").Value = 1000
$ws.Range("AC15").Value = 1000
$ws.Range("AD15").Value = 1000
$ws.Range("AE15").Value = 1000
$ws.Range("AF15").Value = 1000
$ws.Range("AG15").Value = 1000
$ws.Range("AH15").Value = 1000
$ws.Range("AI15").Value = 1000
$ws.Range("AJ15").Value = 1000
$ws.Range("AK15").Value = 1000
$ws.Range("AL15").Value = 1000
$ws.Range("AM15").Value = 1000
$ws.Range("AN15").Value = 1000
$ws.Range("AO15").Value = 1000
